$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserData")

# --- Update the user records shown on the UserData sheet -------------------
# First names first (so the shared-string table gets "Tom3"/"john4" inserted
# before the two new e-mail addresses), then the e-mail addresses.
$ws.Range("A2").Value = "Tom3"
$ws.Range("A3").Value = "john4"
$ws.Range("D2").Value = "tomhanks3@zmail.com"
$ws.Range("D3").Value = "johnswam4@zmail.com"

# --- Widen column D so the longer e-mail addresses are readable ------------
$ws.Columns.Item(4).ColumnWidth = 45.88

# --- Add the new (empty) "GetUserDetails" sheet after "UserData" -----------
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "GetUserDetails"

# --- Restore UserData as the active sheet with the new selected cell -------
$ws.Activate()
$ws.Range("D9").Select()
